$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1:J1 - copy the header style (bold, bordered, centered) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J, rows 2-6
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 2
